$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 82217
$ws.Range("B2").Value = "Sabrina Souza"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45079
$ws.Range("G2").Value = 11664.98

# Row 3
$ws.Range("A3").Value = 9007
$ws.Range("B3").Value = "Milena Nascimento"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45096
$ws.Range("G3").Value = 5161.7

# Row 4
$ws.Range("A4").Value = 13266
$ws.Range("B4").Value = "Ana Lívia Santos"
$ws.Range("C4").Value = "TI"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45101
$ws.Range("G4").Value = 6065.67

# Row 5
$ws.Range("A5").Value = 81426
$ws.Range("B5").Value = "Fernanda Cardoso"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45102
$ws.Range("G5").Value = 5703.76

# Row 6
$ws.Range("A6").Value = 39744
$ws.Range("B6").Value = "Luana Mendes"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45092
$ws.Range("G6").Value = 8795.42

# Row 7
$ws.Range("A7").Value = 92288
$ws.Range("B7").Value = "Dr. Pedro Henrique Almeida"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45101
$ws.Range("G7").Value = 7637.57

# Row 8
$ws.Range("A8").Value = 49580
$ws.Range("B8").Value = "André Silva"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 5461.64

# Row 9
$ws.Range("A9").Value = 57562
$ws.Range("B9").Value = "Elisa Lopes"
$ws.Range("C9").Value = "Operações"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 3903.44

# Row 10
$ws.Range("A10").Value = 76650
$ws.Range("B10").Value = "Miguel Monteiro"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45106
$ws.Range("G10").Value = 11302.8

# Row 11
$ws.Range("A11").Value = 65791
$ws.Range("B11").Value = "Rebeca Mendes"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45097
$ws.Range("G11").Value = 3459.52
